$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C110").Value = 16869
$ws.Range("E110").Value = 25930135
$ws.Range("C115").Value = 17549
$ws.Range("E115").Value = 38605108
$ws.Range("C117").Value = 19702
$ws.Range("E117").Value = 56430548
$ws.Range("C121").Value = 5963
$ws.Range("E121").Value = 11516074
$ws.Range("C122").Value = 9693
$ws.Range("E122").Value = 31946051
$ws.Range("C134").Value = 5673
$ws.Range("E134").Value = 17148182
$ws.Range("C138").Value = 2838
$ws.Range("E138").Value = 6585048
$ws.Range("C139").Value = 3317
$ws.Range("E139").Value = 9227172
$ws.Range("C164").Value = 50565
$ws.Range("E164").Value = 168372196
$ws.Range("C168").Value = 284920
$ws.Range("E168").Value = 1208341610
$ws.Range("C169").Value = 562561
$ws.Range("E169").Value = 1284323670
$ws.Range("C170").Value = 367272
$ws.Range("E170").Value = 2844005381
$ws.Range("C171").Value = 115103
$ws.Range("E171").Value = 444720413
$ws.Range("C173").Value = 54383
$ws.Range("E173").Value = 151846897
$ws.Range("C174").Value = 357160
$ws.Range("E174").Value = 1016615252
$ws.Range("C175").Value = 125508
$ws.Range("E175").Value = 811570009
$ws.Range("C177").Value = 96746
$ws.Range("E177").Value = 174705213
$ws.Range("C179").Value = 235655
$ws.Range("E179").Value = 812084503
$ws.Range("C186").Value = 21933
$ws.Range("E186").Value = 40056540
$ws.Range("C188").Value = 19703
$ws.Range("E188").Value = 66031872
$ws.Range("C196").Value = 7400
$ws.Range("E196").Value = 20652188
$ws.Range("C198").Value = 4508
$ws.Range("E198").Value = 5999357
$ws.Range("C199").Value = 4156
$ws.Range("E199").Value = 9036256
$ws.Range("C203").Value = 13101
$ws.Range("E203").Value = 32996553
$ws.Range("C204").Value = 4754
$ws.Range("E204").Value = 11644170
$ws.Range("C205").Value = 11123
$ws.Range("E205").Value = 44077924
$ws.Range("C209").Value = 5363
$ws.Range("E209").Value = 12210515
$ws.Range("C211").Value = 2863
$ws.Range("E211").Value = 4380689
$ws.Range("C213").Value = 3633
$ws.Range("E213").Value = 11097853
$ws.Range("C214").Value = 6172
$ws.Range("E214").Value = 11075372
$ws.Range("C241").Value = 2583
$ws.Range("E241").Value = 7741099
$ws.Range("C267").Value = 84974
$ws.Range("E267").Value = 156518758
$ws.Range("C295").Value = 91332
$ws.Range("E295").Value = 552911368
$ws.Range("C317").Value = 103579
$ws.Range("E317").Value = 303078087
$ws.Range("C320").Value = 67241
$ws.Range("E320").Value = 124554315
$ws.Range("C322").Value = 81161
$ws.Range("D322").Value = 9703
$ws.Range("E322").Value = 254527631
